# Derive skills from the PR - add new diary entries and clear stray empty cell

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the stray empty string in C23 (was <c r="C23" t="str"><v/></c>)
$ws.Range("C23").ClearContents()

# New diary rows - ensure date column stays plain text (matches existing rows), not auto-converted to a date serial
$ws.Range("A24:A28").NumberFormat = "@"

$ws.Range("A24").Value = "2024-02-12"
$ws.Range("B24").Value = "Created: Fix: Spacing discrepancy on share report modal"

$ws.Range("A25").Value = "2024-02-16"
$ws.Range("B25").Value = "Worked on: Feat: Open Weather Integration"

$ws.Range("A26").Value = "2024-02-19"
$ws.Range("B26").Value = "Worked on: Fix: Action filters default enabled`nMerged: Fix: Action filters default enabled"

$ws.Range("A27").Value = "2024-02-20"
$ws.Range("B27").Value = "Worked on: Feat: Visual password feedback on password forget screen, Feat: Custom field activities and create activites on update/create`nMerged: Feat: Visual password feedback on password forget screen"

$ws.Range("A28").Value = "2024-02-21"
$ws.Range("B28").Value = "Worked on: Feat: Custom field activities and create activites on update/create`nCreated: Fix: Additional email sending on inspection close"
$ws.Range("C28").Value = "Improved my PHP, PHPUnit, Typescript, PHPSpec, and Laravel skills"
